$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:A").Insert()

$ws.Range("A2").Value = "Code Article"
$ws.Range("A3").Value = "E-COM12"
$ws.Range("A4").Value = "FURN_8900"
